# [Poster] Final changes: removed gradient from right
#
# Replace the blue gradient fill (394987 -> 3A4781/3B4785 -> 558BD2) on the
# three "Text Placeholder 19" date/label shapes with a flat solid fill using
# the gradient's final stop color (558BD2).
#
# This COM host doesn't expose a working Fill.Solid()/Fill.Type mutator for
# gradient fills (it only edits the first gradient stop color in that case),
# but assigning ForeColor.RGB to a shape whose fill has just been hidden
# (Fill.Visible = $false, which collapses the fill to <a:noFill/>) correctly
# promotes it to a real <a:solidFill>. That two-step sequence is used below.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 0x558BD2 as a VBA-style RGB() long (R + G*256 + B*65536)
$targetRGB = 13798229

function Set-SolidBlue($shape) {
    $shape.Fill.Visible = $false
    $shape.Fill.ForeColor.RGB = $targetRGB
}

# Shape id 31 "Text Placeholder 19" -> top-level Shapes.Item(11) "Group 37" (id 28)
$grp28 = $s.Shapes.Item(11)
Set-SolidBlue $grp28.GroupItems.Item(2)

# Shape id 37 "Text Placeholder 19" -> top-level Shapes.Item(12) "Group 37" (id 34)
$grp34 = $s.Shapes.Item(12)
Set-SolidBlue $grp34.GroupItems.Item(2)

# Shape id 61 "Text Placeholder 19" -> top-level Shapes.Item(18) "Group 4" (id 5)
$grp5 = $s.Shapes.Item(18)
Set-SolidBlue $grp5.GroupItems.Item(2)
